# Fill in the "保險" (insurance) sheet with the standard metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the other property sheets already carry.
#
# This revision also relabels the "具有相當價值之財產" sheet's bond entries
# from the old "otherbonds" category token to "antique" (the shared-string
# slot the two bond rows pointed at is reused/retextd by this commit).

$wb = $excel.ActiveWorkbook

$wsAntiques = $wb.Worksheets.Item("具有相當價值之財產")
$wsAntiques.Range("F2").Value = "antique"
$wsAntiques.Range("F3").Value = "antique"

$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1) -----------------------------------------------
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# Apply the same bold/bordered header style used on B1:D1 to the new cells
$ws.Range("B1:D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2 --------------------------------------------------------------
$ws.Range("B2").Value = "中華郵政"
$ws.Range("C2").Value = "六六金順"
$ws.Range("D2").Value = "曾巨威"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("H2").Value = "曾巨威"
$ws.Range("I2").Value = 1755
$ws.Range("J2").Value = "tmp8b7f1"
$ws.Range("K2").Value = 158

# --- Row 3 --------------------------------------------------------------
$ws.Range("B3").Value = "富邦人壽"
$ws.Range("C3").Value = "鑫添財萬能終身險"
$ws.Range("D3").Value = "傅寄萍"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("H3").Value = "曾巨威"
$ws.Range("I3").Value = 1755
$ws.Range("J3").Value = "tmp8b7f1"
$ws.Range("K3").Value = 159

# --- Row 4 --------------------------------------------------------------
$ws.Range("B4").Value = "國泰人壽"
$ws.Range("C4").Value = "利率變動型年金(甲型）"
$ws.Range("D4").Value = "傅寄萍"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("H4").Value = "曾巨威"
$ws.Range("I4").Value = 1755
$ws.Range("J4").Value = "tmp8b7f1"
$ws.Range("K4").Value = 160

# Apply the same plain style used on B2:D4 to the new data cells E2:K4
$ws.Range("B2:D4").Copy()
$ws.Range("E2:K4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The "date" column holds a plain text token ("2012-04-16"), not a real
# date value, so force text formatting before assigning, which stops
# Excel's autodetection from rewriting it as a date serial number.
$ws.Range("G2:G4").NumberFormat = "@"
$ws.Range("G2").Value = "2012-04-16"
$ws.Range("G3").Value = "2012-04-16"
$ws.Range("G4").Value = "2012-04-16"
